# Apply the documented text edits to the gendoc template.
$d = $word.ActiveDocument

# 1) Main body: the default output-file timestamp suffix is replaced with a
#    literal "1" so the generated file is named "...+gendoc.1.docx" instead
#    of using the ${date}/${time} placeholders.
$d.Content.Find.Execute(
    "`${date}.`${time}docx'", $true, $false, $false, $false, $false,
    $true, 1, $false, "1.docx'", 2
) | Out-Null

# 2) Header: bump the trailing version-suffix placeholder from ".n" to ".1"
#    so the header shows "...+gendoc.1" (matching the new profile version).
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hdr.Range.Find.Execute(
        "gendoc.n", $true, $false, $false, $false, $false,
        $true, 1, $false, "gendoc.1", 2
    ) | Out-Null
}
